$wb = $excel.ActiveWorkbook

# The new sheet should be placed right after the existing "Login" sheet.
$loginSheet = $wb.Worksheets.Item("Login")

# Add the new "Urls" worksheet.
$newSheet = $wb.Worksheets.Add($null, $loginSheet)
$newSheet.Name = "Urls"

# Populate the lookup table of page names -> URLs.
$newSheet.Range("A1").Value = "Login"
$newSheet.Range("B1").Value = "https://demo.actitime.com/login.do"
$newSheet.Range("A2").Value = "Dashboard"
$newSheet.Range("B2").Value = "https://demo.actitime.com/user/submit_tt.do"

# Make the new "Urls" sheet the active tab / selected sheet.
$newSheet.Select() | Out-Null
$newSheet.Range("K17").Select() | Out-Null

$wb.Save()
